# PITCH.pptx - "Correção do PITCH versão final - ajuste do LOGO"
#
# The logo placeholder shape on slide 1 is converted from an ellipse
# ("Oval 3") into a rectangle ("Rectangle 3"), repositioned/resized, and
# its caption text is trimmed from "LOGO DO PROJETO" down to just "LOGO"
# (leaving a trailing blank line where the second line of text used to
# continue).

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# The logo shape is the 3rd shape on the title slide (id=4, "Oval 3").
$logo = $s.Shapes.Item(3)

# Rename it to match its new look.
$logo.Name = "Rectangle 3"

# Reposition / resize (values chosen so the point->EMU conversion lands
# exactly on 1354666,3522134 / 2032001,1219199 EMU).
$logo.Left = 106.66662
$logo.Top = 277.3334
$logo.Width = 160.0001
$logo.Height = 95.99992

# Swap the oval geometry for a plain rectangle.
$logo.AutoShapeType = 1   # msoShapeRectangle (was 9, msoShapeOval)

# Shorten the caption; press-enter leaves a trailing empty paragraph.
$logo.TextFrame.TextRange.Text = "LOGO`r"

# --- Best-effort cleanup -----------------------------------------------
# The committed version also drops the (unused) notes master part - this
# presentation's only slide has no real speaker notes, so PowerPoint
# considers the notes master orphaned and removes both
# ppt/notesMasters/notesMaster1.xml and its private theme
# (ppt/theme/theme2.xml), along with the <p:notesMasterIdLst> entry in
# presentation.xml. Attempt it where the host supports it; ignore
# failures so the main edit above is unaffected on hosts that don't.
try {
    $p.NotesMaster.Delete()
} catch {
}
